$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G6").NumberFormat = "0.00%"
$ws.Range("G9").NumberFormat = "0.00%"

$ws.Range("G2").Value = 29.3/100
$ws.Range("G3").Value = 9.3/100
$ws.Range("G4").Value = 89.7/100
$ws.Range("G5").Value = 0/100
$ws.Range("G6").Value = 41.5/100
$ws.Range("G7").Value = 1.3
$ws.Range("G9").Value = 66.3/100

$ws.Columns("G").AutoFit()

$ws.Range("G9").Select()
